# Update "Pais" sheet with refreshed COVID-19 country data (10 Aug 2020, 16:25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Agosto de 2020 a las 16:25"

# --- Country name re-ordering (rows whose rank/name changed) ---
# Namibia / Mayotte / Cuba block
$ws.Range("A118").Value = "Namibia"
$ws.Range("A119").Value = "Mayotte"
$ws.Range("A120").Value = "Cuba"

# Timor Oriental / Santa Lucia swap
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"

# Montserrat / Islas Malvinas swap
$ws.Range("A213").Value = "Montserrat"
$ws.Range("A214").Value = "Islas Malvinas"

# --- Numeric data refresh (B:Casos totales, C:Nuevos casos, D:Casos activos,
#     E:Recuperados, F:Casos criticos, G:Muertes hoy, H:Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 5201933
$ws.Range("C4").Value = 2489
$ws.Range("D4").Value = 2665033
$ws.Range("E4").Value = 2371272
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 165628

# Row 20 - Argentina
$ws.Range("E20").Value = 133623
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = 4634

# Row 22 - Alemania
$ws.Range("B22").Value = 217306
$ws.Range("C22").Value = 25
$ws.Range("E22").Value = 10146
$ws.Range("H22").Value = 9260

# Row 28 - Catar
$ws.Range("B28").Value = 113262
$ws.Range("C28").Value = 315
$ws.Range("D28").Value = 109993
$ws.Range("E28").Value = 3081
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 188

# Row 35 - Suecia
$ws.Range("B35").Value = 82972
$ws.Range("G35").Value = 5
$ws.Range("H35").Value = 5766

# Row 48 - Portugal
$ws.Range("B48").Value = 52825
$ws.Range("C48").Value = 157
$ws.Range("D48").Value = 38600
$ws.Range("E48").Value = 12466
$ws.Range("G48").Value = 3
$ws.Range("H48").Value = 1759

# Row 60 - Azerbaiyan
$ws.Range("B60").Value = 33647
$ws.Range("C60").Value = 79
$ws.Range("D60").Value = 30642
$ws.Range("E60").Value = 2513
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 492

# Row 62 - Uzbekistan
$ws.Range("B62").Value = 31068
$ws.Range("C62").Value = 459
$ws.Range("D62").Value = 22559
$ws.Range("E62").Value = 8311
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 198

# Row 63 - Serbia
$ws.Range("B63").Value = 28262
$ws.Range("C63").Value = 163
$ws.Range("E63").Value = 13569
$ws.Range("G63").Value = 5
$ws.Range("H63").Value = 646

# Row 86 - Noruega
$ws.Range("B86").Value = 9661
$ws.Range("C86").Value = 23
$ws.Range("E86").Value = 548

# Row 93 - Tayikistan
$ws.Range("B93").Value = 7827
$ws.Range("C93").Value = 82
$ws.Range("D93").Value = 6614
$ws.Range("E93").Value = 1151

# Row 118 - Namibia (values after the name/rank update above)
$ws.Range("B118").Value = 3101
$ws.Range("C118").Value = 152
$ws.Range("D118").Value = 715
$ws.Range("E118").Value = 2367
$ws.Range("H118").Value = 19

# Row 119 - Mayotte
$ws.Range("B119").Value = 3068
$ws.Range("D119").Value = 2835
$ws.Range("E119").Value = 194
$ws.Range("H119").Value = 39

# Row 120 - Cuba
$ws.Range("B120").Value = 2953
$ws.Range("D120").Value = 2451
$ws.Range("E120").Value = 414
$ws.Range("H120").Value = 88

# Row 154 - Principado de Andorra
$ws.Range("B154").Value = 963
$ws.Range("C154").Value = 8
$ws.Range("E154").Value = 72

# Row 178 - Trinidad yTobago
$ws.Range("B178").Value = 280
$ws.Range("C178").Value = 1
$ws.Range("D178").Value = 138
$ws.Range("E178").Value = 134

# Row 213 - Montserrat
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

# Row 214 - Islas Malvinas
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
